$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E").EntireColumn.Insert()
$ws.Range("E2").Borders.LineStyle = -4142
Write-Output "done"
